# Fruta / hortaliza, semanal
# The rows 2-32 (one record each) get reshuffled: the values in columns
# D (Fecha), J (Volumen), K (Precio mínimo), L (Precio máximo),
# M (Precio promedio ponderado) and P (Precio $/Kg) move to a different
# row while the rest of the row (Mercado, Región, Categoría, etc.) stays
# put. Capture the "before" values for those columns first, then write
# them back out per the permutation map (row -> source row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> source row (which row's D/J/K/L/M/P values land here)
$map = @{
    2 = 24; 3 = 22; 4 = 20; 5 = 6; 6 = 21; 7 = 17; 8 = 31; 9 = 9; 10 = 15;
    11 = 3; 12 = 2; 13 = 26; 14 = 10; 15 = 8; 16 = 12; 17 = 25; 18 = 32;
    19 = 5; 20 = 11; 21 = 30; 22 = 18; 23 = 29; 24 = 16; 25 = 7; 26 = 28;
    27 = 14; 28 = 4; 29 = 13; 30 = 23; 31 = 19; 32 = 27
}

$cols = @(4, 10, 11, 12, 13, 16)   # D, J, K, L, M, P

# Snapshot the original values before any writes.
$orig = @{}
for ($r = 2; $r -le 32; $r++) {
    $row = @{}
    foreach ($c in $cols) {
        $row[$c] = $ws.Cells.Item($r, $c).Value2
    }
    $orig[$r] = $row
}

# Write back according to the permutation map.
for ($r = 2; $r -le 32; $r++) {
    $src = $map[$r]
    $srcRow = $orig[$src]
    foreach ($c in $cols) {
        $ws.Cells.Item($r, $c).Value = $srcRow[$c]
    }
}
